# Swap the presentation's colour theme from the "Integral" (Red Violet)
# scheme over to the stock "Office Theme" palette, per the commit:
#   theme1.xml (the deck's real, slide-master theme) goes from the
#   "Red Violet" colours to the default "Office" colours; what used to
#   sit in theme2.xml (the Office Theme palette) becomes the live theme.
#
# PowerPoint's object model doesn't expose the theme/colour-scheme
# display *names* for editing (Theme.Name and the <a:clrScheme name=.../>
# bookkeeping string are not writable through COM - they're only ever
# stamped by the UI when a built-in theme is chosen from the gallery),
# so this reproduces the edit the way real COM automation would: by
# pushing the twelve standard theme colour slots to the "Office" values
# via ColorScheme.Colors(n).RGB, using the classic COM RGB(r,g,b) packing
# (R + G*256 + B*65536).

function ComRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

# index -> (slot name, target "Office" hex colour)
$scheme.Colors(1).RGB  = ComRGB 0x00 0x00 0x00   # dk1      -> 000000
$scheme.Colors(2).RGB  = ComRGB 0xFF 0xFF 0xFF   # lt1      -> FFFFFF
$scheme.Colors(3).RGB  = ComRGB 0x44 0x54 0x6A   # dk2      -> 44546A
$scheme.Colors(4).RGB  = ComRGB 0xE7 0xE6 0xE6   # lt2      -> E7E6E6
$scheme.Colors(5).RGB  = ComRGB 0x5B 0x9B 0xD5   # accent1  -> 5B9BD5
$scheme.Colors(6).RGB  = ComRGB 0xED 0x7D 0x31   # accent2  -> ED7D31
$scheme.Colors(7).RGB  = ComRGB 0xA5 0xA5 0xA5   # accent3  -> A5A5A5
$scheme.Colors(8).RGB  = ComRGB 0xFF 0xC0 0x00   # accent4  -> FFC000
$scheme.Colors(9).RGB  = ComRGB 0x44 0x72 0xC4   # accent5  -> 4472C4
$scheme.Colors(10).RGB = ComRGB 0x70 0xAD 0x47   # accent6  -> 70AD47
$scheme.Colors(11).RGB = ComRGB 0x05 0x63 0xC1   # hlink    -> 0563C1
$scheme.Colors(12).RGB = ComRGB 0x95 0x4F 0x72   # folHlink -> 954F72
